# Regenerate merged AHB files
# Rename the header row from the generic "_old"/"_new" suffixes to the
# concrete version tags ("_FV2210" / "_FV2304"), then (re)build the table
# over the data range and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

# Columns A-J: "<name>_old" -> "<name>_FV2210"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $baseNames[$i] + "_FV2210"
}

# Column K stays "diff" (unchanged)

# Columns L-U: "<name>_new" -> "<name>_FV2304"
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Cells.Item(1, 11 + $i + 1).Value = $baseNames[$i] + "_FV2304"
}

# Turn the data range into a real Excel Table (ListObject) with an AutoFilter.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U74"), $false, 1, $null)
$tbl.Name = "Table1"

# Freeze the header row.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
